$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The averaged-intensities generation script was re-run including a
# "Gaussian-Quadrature" scheme (now correctly grouped with the other
# special/ring schemes) plus three brand-new spiral sampling schemes, which
# pushes the previously-trailing schemes (NoRotation/Rotation/HexGrid...)
# further down the table and appends the HexGrid rows at the bottom.

# Re-label existing rows 10-16 (column B) to reflect the new scheme order.
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"

# Append the three trailing rows (17-19), re-surfacing the schemes that used
# to sit at rows 13-15 (HexGrid-*), each with an averaged intensity of 1
# across all 11 HKL columns (C:M).
$newRows = @(
    @{ Row = 17; A = 15; B = "HexGrid-90degTilt5degRes" },
    @{ Row = 18; A = 16; B = "HexGrid-90degTilt22p5degRes" },
    @{ Row = 19; A = 17; B = "HexGrid-60degTilt5degRes" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A keeps the same bold/centered/bordered look as the rows above -
    # copy formatting from the last existing row (A16) then set the value.
    $ws.Range("A16").Copy() | Out-Null
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Cells.Item($row, 2).Value = $r.B
    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}

$excel.CutCopyMode = 0
